# Update the active selection shown in the sheet view (N12 -> B21)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray prior value that had been entered in B14
$ws.Range("B14").ClearContents()

# "Go fish" prior probabilities used for the entropy calc at the bottom of
# the sheet (row 19) - update counts/percentages (formula in E19 recalculates
# automatically from these inputs)
$ws.Range("B19").Value = 0.4
$ws.Range("C19").Value = 0.4
$ws.Range("D19").Value = 0.2

# Move/restore the active cell selection to B21
$ws.Range("B21").Select()
